$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IP")

# Row 4 - Inventory
$ws.Range("B4").Value = 1828000000.0
$ws.Range("C4").Value = 2050000000.0
$ws.Range("D4").Value = 2007000000.0
$ws.Range("E4").Value = 2010000000.0
$ws.Range("F4").Value = 2002000000.0

# Row 14 - Accounts Payable
$ws.Range("B14").Value = 3805000000.0
$ws.Range("C14").Value = 2320000000.0
$ws.Range("D14").Value = 2226000000.0
$ws.Range("E14").Value = 2206000000.0
$ws.Range("F14").Value = 2379000000.0

# Row 15 - Accrued Expenses (was an empty inline string, now a number)
$ws.Range("B15").Value = 363000000.0

# Row 21 - Long Term Tax Liability (Deferred)
$ws.Range("B21").Value = 2756000000.0
$ws.Range("C21").Value = 2743000000.0
$ws.Range("D21").Value = 2639000000.0
$ws.Range("E21").Value = 2654000000.0
$ws.Range("F21").Value = 2662000000.0

# Row 26 - Additional Paid In Capital (was an empty inline string, now a number)
$ws.Range("B26").Value = 6267000000.0

# Row 29 - Treasury Stock (was an empty inline string, now a number)
$ws.Range("B29").Value = 2719000000.0

# Row 34 - Shares (Common) (was an empty inline string, now a number)
$ws.Range("B34").Value = 391700000.0

# Row 36 - Net Debt
$ws.Range("G36").Value = 13608000000.0

# Row 37 - Total Debt
$ws.Range("G37").Value = 14119000000.0
